$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Translate the header row from Spanish to English
$ws.Range("A1").Value = "Language"
$ws.Range("B1").Value = "Speakers"

# Update the ListObject (Table1) column headers to match
$table = $ws.ListObjects.Item(1)
$table.ListColumns.Item(1).Name = "Language"
$table.ListColumns.Item(2).Name = "Speakers"

# Narrow column B slightly (target stored width ~10.71 chars; ColumnWidth is
# quantized to the sheet's pixel grid, so 9.8 is the closest input that lands
# on the nearest achievable width)
$ws.Columns.Item(2).ColumnWidth = 9.8
